$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.456.07"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.107.61"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D5").Value = "'336.55"
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "'0.5245"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("D8").Value = "'0.4590"
$ws.Range("E8").Value = "  +5.73%  "
$ws.Range("D9").Value = "'53.37"
$ws.Range("E9").Value = "  +15.57%  "
$ws.Range("D10").Value = "'0.08953"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("D12").Value = "'24.48"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "2.092.47"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").Value = "'6.804"
$ws.Range("E14").Value = "  +2.10%  "
$ws.Range("D15").Value = "'7.966"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").Value = "'96.60"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "'0.00001133"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").Value = "'0.06637"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "'19.32"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").Value = "'6.308"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").Value = "30.520.93"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").Value = "'12.37"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "'2.365"
$ws.Range("E25").Value = "  +2.85%  "
$ws.Range("D26").Value = "2.339.94"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "'22.36"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D28").Value = "'2.572"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").Value = "'163.67"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").Value = "'132.85"
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("D31").Value = "'1.217"
$ws.Range("E31").Value = "  +2.32%  "
$ws.Range("D34").Value = "'6.203"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("D35").Value = "'3.924"
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("D36").Value = "'10.48"
$ws.Range("E36").Value = "  +8.28%  "
$ws.Range("D37").Value = "'0.02579"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").Value = "'0.06844"
$ws.Range("E38").Value = "  +2.68%  "
$ws.Range("D39").Value = "'5.561"
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("D40").Value = "'12.87"
$ws.Range("E40").Value = "  +3.29%  "
$ws.Range("D41").Value = "'0.2301"
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("D42").Value = "'0.6903"
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("D43").Value = "'1.246"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Value = "'2.353"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D48").Value = "'3.668"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").Value = "'0.00000000354"
$ws.Range("E49").Value = "  +25.43%  "
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "'83.67"
$ws.Range("E51").Value = "  +1.71%  "

# Rows 32/33 and 46/47 content swapped (new coins inserted at same rank position)
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1075"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'1.702"
$ws.Range("E33").Value = "  +10.91%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'14.09"
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.6392"
$ws.Range("E47").Value = "  +0.84%  "
